$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.438.42"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.602.51"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.98"
$ws.Range("E5").Value = "  -2.37%  "

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.40"
$ws.Range("E6").Value = "  -3.22%  "

# Row 7: XRP -> XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  -2.30%  "

# Row 8: LidoStakedEther -> LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.595.90"
$ws.Range("E8").Value = "  -0.87%  "

# Row 9: USDC -> USDC
$ws.Range("E9").Value = "  +0.12%  "

# Row 10: Dogecoin -> Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  +4.11%  "

# Row 11: Cardano -> Cardano
$ws.Range("E11").Value = "  -1.35%  "

# Row 12: Avalanche -> Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.94"
$ws.Range("E12").Value = "  -4.53%  "

# Row 13: ShibaInu -> ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000312"
$ws.Range("E13").Value = "  +7.19%  "

# Row 14: Polkadot -> Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.68"
$ws.Range("E14").Value = "  -2.69%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.177.59"
$ws.Range("E15").Value = "  -0.82%  "

# Row 16: Chainlink -> Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.81"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17: WrappedEther -> WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.588.50"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18: WrappedBTC -> WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.366.94"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19: Uniswap -> Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  -0.58%  "

# Row 20: TRON -> TRON
$ws.Range("E20").Value = "  +0.07%  "

# Row 21: Polygon -> Polygon
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("E21").Value = "  -1.93%  "

# Row 22: BitcoinCash -> BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.74"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.47"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24: Toncoin -> Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.87"
$ws.Range("E24").Value = "  -9.55%  "

# Row 25: Litecoin -> Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.98"
$ws.Range("E25").Value = "  +6.19%  "

# Row 26: PancakeSwap -> PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.36"
$ws.Range("E26").Value = "  -2.49%  "

# Row 27: ImmutableX -> ImmutableX
$ws.Range("E27").Value = "  -5.88%  "

# Row 28: RenderToken -> RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.02"
$ws.Range("E28").Value = "  -4.88%  "

# Row 29: Filecoin -> Filecoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("E29").Value = "  -3.25%  "

# Row 30: EthereumClassic -> EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.28"
$ws.Range("E30").Value = "  -2.22%  "

# Row 31: NEARProtocol -> NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.62"
$ws.Range("E31").Value = "  -4.06%  "

# Row 32: Cosmos -> Cosmos
$ws.Range("E32").Value = "  -0.59%  "

# Row 33: Hedera -> Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("E33").Value = "  -2.79%  "

# Row 34: OKB -> OKB
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.76"
$ws.Range("E34").Value = "  -0.46%  "

# Row 35: Bittensor -> Bittensor
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "575.57"
$ws.Range("E35").Value = "  -8.29%  "

# Row 36: InjectiveProtocol -> InjectiveProtocol
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.80"
$ws.Range("E36").Value = "  -4.54%  "

# Row 37: Dai -> PEPE
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0813"
$ws.Range("E37").Value = "  -1.54%  "

# Row 38: PEPE -> Dai
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.16%  "

# Row 39: TheGraph -> TheGraph
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.397"
$ws.Range("E39").Value = "  -4.57%  "

# Row 40: Fetch.AI -> dogwifhat
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  +15.77%  "

# Row 41: dogwifhat -> Fetch.AI
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").Value = "  +4.60%  "

# Row 42: Stacks -> Stacks
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  -3.22%  "

# Row 43: Kaspa -> Kaspa
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.138"
$ws.Range("E43").Value = "  -6.30%  "

# Row 44: ThetaToken -> ThetaToken
$ws.Range("E44").Value = "  -3.73%  "

# Row 45: Maker -> ApeXProtocol
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.54"
$ws.Range("E45").Value = "  +7.21%  "

# Row 46: VeChain -> Maker
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.215.25"
$ws.Range("E46").Value = "  -2.36%  "

# Row 47: ApeXProtocol -> VeChain
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0445"
$ws.Range("E47").Value = "  -1.93%  "

# Row 48: THORChain -> THORChain
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.68"
$ws.Range("E48").Value = "  +4.74%  "

# Row 49: Stellar -> Stellar
$ws.Range("E49").Value = "  -0.63%  "

# Row 50: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("E50").Value = "  -0.08%  "

# Row 51: LidoDAOToken -> LidoDAOToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.23"
$ws.Range("E51").Value = "  -3.77%  "
